$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds the "Location OK/KO" header (row 1) and "KO" values (rows 2:2637).
# Commit: "Changing geotraining OK/KO to yes/no"
$colD = $ws.Range("D1:D2637")

# 1. Rename the header cell D1: "Location OK/KO" -> "Location yes/no"
$colD.Replace("Location OK/KO", "Location yes/no", [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole)

# 2. Replace every remaining whole-cell "KO" value in the column with "no"
$colD.Replace("KO", "no", [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole)

# 3. Reflect the saved view state: active cell/selection moved from E1 to D1
$ws.Range("D1").Select()

# 4. The header row's auto-fit height shrank slightly (35.8 -> 35.05) once the
#    shorter "Location yes/no" text replaced "Location OK/KO" in the wrapped cell
$ws.Rows.Item(1).RowHeight = 35.05

Write-Host "Updated geolocation header and values (OK/KO -> yes/no)"
